# Generate Report for Handback
# Replaces the stale handback identifiers/timestamps in the report with the
# freshly generated ones (new GUID-named files + new xlf names + new dates).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: update the TextToDisplay of the hyperlink anchored at a given cell
# ---------------------------------------------------------------------------
# NOTE: this COM-interop PowerShell engine does not reliably bind named
# (-Param value) arguments on user functions, so this helper (and all call
# sites below) use positional parameters only: Sheet, CellAddress, NewDisplay.
function Set-HyperlinkDisplay {
    param($Sheet, $CellAddress, $NewDisplay)
    $target = $CellAddress
    if ($CellAddress -match '^([A-Za-z]+)(\d+)$') {
        $col = $matches[1]
        $row = $matches[2]
        $target = "`$$col`$$row"
    }
    foreach ($hl in $Sheet.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $target) {
            $hl.TextToDisplay = $NewDisplay
        }
    }
}

# Old identifiers being replaced
$oldGuid1 = "4feb7aad-90be-416b-9afe-7194f5e88201"
$oldGuid2 = "b003a3cc-756c-4406-ab57-bbec474a61c0"

# New identifiers
$newGuid1 = "e1c06d56-8123-44e7-8b18-d208dc87f51c"
$newGuid2 = "ffff2e350816-22a7-49af-af69-054f76310a7c"

$oldZhXlf = "$oldGuid2.84b3f4a4a727a2e5672184347ddb206731031e31.zh-cn.xlf"
$oldDeXlf = "$oldGuid1.1bc369c1e67c487b37f3400ca5226a774266676f.de-de.xlf"

$newZhXlf = "$newGuid1.00fdb81ad85f88e42d024a98e1866d0cd957f91d.zh-cn.xlf"
$newDeXlf = "$newGuid1.00fdb81ad85f88e42d024a98e1866d0cd957f91d.de-de.xlf"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
Set-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-09-01 01:06:49"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
Set-HyperlinkDisplay $wsOverview "B3" "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-09-01 01:06:49"

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
Set-HyperlinkDisplay $wsZhCn "A2" "$newGuid1.md"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
Set-HyperlinkDisplay $wsZhCn "I2" "$newGuid1.md"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
Set-HyperlinkDisplay $wsZhCn "A3" "$newGuid2.md"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
Set-HyperlinkDisplay $wsZhCn "I3" "$newGuid2.md"

$wsZhCn.Range("G3").Value = $newZhXlf
$wsZhCn.Range("J3").Value = $newZhXlf

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
Set-HyperlinkDisplay $wsDeDe "A2" "$newGuid1.md"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
Set-HyperlinkDisplay $wsDeDe "I2" "$newGuid1.md"

$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-09-01 01:06:49"
$wsDeDe.Range("J2").Value = $newDeXlf
$wsDeDe.Range("K2").Value = "2016-09-01 01:07:16"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
Set-HyperlinkDisplay $wsDeDe "A3" "$newGuid2.md"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
Set-HyperlinkDisplay $wsDeDe "I3" "$newGuid2.md"

$wsDeDe.Range("G3").Value = $newDeXlf
$wsDeDe.Range("H3").Value = "2016-09-01 01:06:49"
$wsDeDe.Range("J3").Value = $newDeXlf
$wsDeDe.Range("K3").Value = "2016-09-01 01:07:16"
